# Fix : 몬스터 Load 함수
# Update monster Name / PrefabPath columns on the MonsterData sheet:
#  - rows 6,7 (Stega type): PachycephalaRed-style duplicate "StegaBlue" names
#    become distinct StegaGreen / StegaPurple
#  - rows 9,10 (Rapto type): duplicate "RaptoBlue" names become distinct
#    RaptoGreen / RaptoOrange
#  - PrefabPath column (H) gets a "Prefabs/" prefix and the per-row color
#    suffix corrected to match the new Name values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Name column (C) fixes --------------------------------------------------
$ws.Range("C6").Value  = "StegaGreen"
$ws.Range("C7").Value  = "StegaPurple"
$ws.Range("C9").Value  = "RaptoGreen"
$ws.Range("C10").Value = "RaptoOrange"

# -- PrefabPath column (H) fixes --------------------------------------------
$ws.Range("H2").Value  = "Prefabs/Monster/1/Blue"
$ws.Range("H3").Value  = "Prefabs/Monster/1/Green"
$ws.Range("H4").Value  = "Prefabs/Monster/1/Red"
$ws.Range("H5").Value  = "Prefabs/Monster/2/Blue"
$ws.Range("H6").Value  = "Prefabs/Monster/2/Green"
$ws.Range("H7").Value  = "Prefabs/Monster/2/Purple"
$ws.Range("H8").Value  = "Prefabs/Monster/3/Blue"
$ws.Range("H9").Value  = "Prefabs/Monster/3/Green"
$ws.Range("H10").Value = "Prefabs/Monster/3/Orange"

# -- View state: selection moved to H8, zoomed in to 160% ------------------
$ws.Range("H8").Select()
$excel.ActiveWindow.Zoom = 160
